$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (price) to text format while writing values that could
# otherwise be auto-detected as numbers (e.g. single-decimal-point values),
# then restore the original (default/general) style so the saved XML has no
# style index on these cells, matching the original formatting.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '26.446.36'
$ws.Range('E2').Value = '  +1.53%  '
$ws.Range('D3').Value = '1.675.69'
$ws.Range('E3').Value = '  +2.39%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '216.94'
$ws.Range('E5').Value = '  +1.67%  '
$ws.Range('E6').Value = '  +1.45%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '0.2692'
$ws.Range('E8').Value = '  +3.68%  '
$ws.Range('D9').Value = '0.06406'
$ws.Range('E9').Value = '  +1.79%  '
$ws.Range('D10').Value = '21.79'
$ws.Range('E10').Value = '  +5.29%  '
$ws.Range('D11').Value = '0.07799'
$ws.Range('E11').Value = '  +1.86%  '
$ws.Range('D12').Value = '1.681.45'
$ws.Range('E12').Value = '  +3.00%  '
$ws.Range('D13').Value = '4.512'
$ws.Range('E13').Value = '  +2.18%  '
$ws.Range('D14').Value = '0.5573'
$ws.Range('E14').Value = '  +0.75%  '
$ws.Range('D15').Value = '0.0₅8318'
$ws.Range('E15').Value = '  +1.04%  '
$ws.Range('D16').Value = '65.65'
$ws.Range('E16').Value = '  +1.06%  '
$ws.Range('D17').Value = '26.504.17'
$ws.Range('E17').Value = '  +1.82%  '
$ws.Range('E18').Value = '  -0.07%  '
$ws.Range('D19').Value = '4.777'
$ws.Range('E19').Value = '  +1.75%  '
$ws.Range('D20').Value = '193.51'
$ws.Range('E20').Value = '  +3.05%  '
$ws.Range('D21').Value = '10.29'
$ws.Range('E21').Value = '  +1.12%  '
$ws.Range('D22').Value = '6.335'
$ws.Range('E22').Value = '  +2.92%  '
$ws.Range('E23').Value = '  +0.05%  '
$ws.Range('D24').Value = '142.42'
$ws.Range('E24').Value = '  -1.93%  '
$ws.Range('D25').Value = '0.1279'
$ws.Range('E25').Value = '  +5.16%  '
$ws.Range('D26').Value = '7.405'
$ws.Range('E26').Value = '  -0.10%  '
$ws.Range('D27').Value = '16.27'
$ws.Range('E27').Value = '  +3.18%  '
$ws.Range('E28').Value = '  +3.93%  '
$ws.Range('D29').Value = '0.06280'
$ws.Range('E29').Value = '  +5.43%  '
$ws.Range('D30').Value = '1.277'
$ws.Range('E30').Value = '  +2.00%  '
$ws.Range('D31').Value = '3.615'
$ws.Range('E31').Value = '  +5.20%  '
$ws.Range('D32').Value = '3.451'
$ws.Range('E32').Value = '  +1.10%  '
$ws.Range('D33').Value = '1.687'
$ws.Range('E33').Value = '  +2.64%  '
$ws.Range('D34').Value = '1.008'
$ws.Range('E34').Value = '  +2.37%  '
$ws.Range('D35').Value = '0.6197'
$ws.Range('E35').Value = '  +9.37%  '
$ws.Range('E36').Value = '  +0.90%  '
$ws.Range('E37').Value = '  +0.90%  '
$ws.Range('D38').Value = '6.190'
$ws.Range('E38').Value = '  +7.45%  '
$ws.Range('D39').Value = '0.01637'
$ws.Range('E39').Value = '  +1.39%  '
$ws.Range('D40').Value = '1.095.01'
$ws.Range('E40').Value = '  +5.77%  '
$ws.Range('D41').Value = '0.8653'
$ws.Range('E41').Value = '  +1.79%  '
$ws.Range('D42').Value = '0.9999'
$ws.Range('E42').Value = '  -0.09%  '
$ws.Range('D43').Value = '100.64'
$ws.Range('E43').Value = '  +0.38%  '
$ws.Range('D44').Value = '1.821.17'
$ws.Range('E44').Value = '  +2.06%  '
$ws.Range('D45').Value = '57.75'
$ws.Range('E45').Value = '  +3.74%  '
$ws.Range('D46').Value = '8.139'
$ws.Range('E46').Value = '  +1.24%  '
$ws.Range('D47').Value = '0.9999'
$ws.Range('E47').Value = '  -0.21%  '
$ws.Range('D48').Value = '0.0₈103'
$ws.Range('E48').Value = '  -4.22%  '
$ws.Range('D49').Value = '0.05211'
$ws.Range('E49').Value = '  +0.98%  '
$ws.Range('D50').Value = '1.482'
$ws.Range('E50').Value = '  +6.91%  '
$ws.Range('D51').Value = '6.048'
$ws.Range('E51').Value = '  +2.21%  '

# Restore column D cells to the default (unstyled) format so no stray
# style attribute gets written to the saved worksheet XML.
$ws.Range("D2:D51").Style = $ws.Range("D4").Style

